# Auto-generated: apply scheduled market-data refresh values to FFXIV leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 151.45454
$ws.Range("I12").Value = 108.71429
$ws.Range("J12").Value = 226.25
$ws.Range("K12").Value = 108.71429
$ws.Range("L12").Value = 226.25
$ws.Range("M12").Value = 61.28570999999999
$ws.Range("N12").Value = -566.25
$ws.Range("H55").Value = 91882.17999999999
$ws.Range("I55").Value = 142958.14
$ws.Range("J55").Value = 2499.25
$ws.Range("K55").Value = 142958.14
$ws.Range("L55").Value = 2499.25
$ws.Range("M55").Value = -142744.14
$ws.Range("N55").Value = -2927.25
$ws.Range("H58").Value = 4280.9487
$ws.Range("J58").Value = 5281.484
$ws.Range("L58").Value = 15844.452
$ws.Range("N58").Value = -16144.452
$ws.Range("H74").Value = 11799.353
$ws.Range("J74").Value = 14374.75
$ws.Range("L74").Value = 14374.75
$ws.Range("N74").Value = -16246.75
$ws.Range("H77").Value = 11799.353
$ws.Range("J77").Value = 14374.75
$ws.Range("L77").Value = 71873.75
$ws.Range("N77").Value = -81233.75
$ws.Range("H80").Value = 922.1053000000001
$ws.Range("I80").Value = 960.3
$ws.Range("J80").Value = 879.6667
$ws.Range("K80").Value = 2880.9
$ws.Range("L80").Value = 2639.0001
$ws.Range("M80").Value = -1882.9
$ws.Range("N80").Value = -4635.0001
$ws.Range("H83").Value = 922.1053000000001
$ws.Range("I83").Value = 960.3
$ws.Range("J83").Value = 879.6667
$ws.Range("K83").Value = 8642.699999999999
$ws.Range("L83").Value = 7917.0003
$ws.Range("M83").Value = -3650.699999999999
$ws.Range("N83").Value = -17901.0003
$ws.Range("H86").Value = 2927319.5
$ws.Range("I86").Value = 2147
$ws.Range("J86").Value = 5267457.5
$ws.Range("K86").Value = 2147
$ws.Range("L86").Value = 5267457.5
$ws.Range("M86").Value = -1024
$ws.Range("N86").Value = -5269703.5
$ws.Range("H89").Value = 2927319.5
$ws.Range("I89").Value = 2147
$ws.Range("J89").Value = 5267457.5
$ws.Range("K89").Value = 10735
$ws.Range("L89").Value = 26337287.5
$ws.Range("M89").Value = -5119
$ws.Range("N89").Value = -26348519.5
$ws.Range("H132").Value = 10242.7295
$ws.Range("I132").Value = 1860.738
$ws.Range("J132").Value = 68916.664
$ws.Range("K132").Value = 5582.214
$ws.Range("L132").Value = 206749.992
$ws.Range("M132").Value = -3052.214
$ws.Range("N132").Value = -211809.992
$ws.Range("H138").Value = 5871.8384
$ws.Range("J138").Value = 6284.259
$ws.Range("L138").Value = 18852.777
$ws.Range("N138").Value = -29132.777

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1988.7
$ws.Range("J88").Value = 791.3077
$ws.Range("L88").Value = 791.3077
$ws.Range("N88").Value = -1603.3077
$ws.Range("H91").Value = 1988.7
$ws.Range("J91").Value = 791.3077
$ws.Range("L91").Value = 791.3077
$ws.Range("N91").Value = -3599.3077
$ws.Range("H132").Value = 3208.46
$ws.Range("I132").Value = 3034.0789
$ws.Range("K132").Value = 9102.236699999999
$ws.Range("M132").Value = -6572.236699999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2033.3529
$ws.Range("I99").Value = 1791.5
$ws.Range("K99").Value = 1791.5
$ws.Range("M99").Value = -293.5
$ws.Range("H105").Value = 2925.875
$ws.Range("I105").Value = 2925.875
$ws.Range("K105").Value = 2925.875
$ws.Range("M105").Value = -1178.875
$ws.Range("H107").Value = 1432224.1
$ws.Range("I107").Value = 2989.5
$ws.Range("K107").Value = 2989.5
$ws.Range("M107").Value = -1069.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4840.25
$ws.Range("I62").Value = 3227.5
$ws.Range("K62").Value = 3227.5
$ws.Range("M62").Value = -2603.5
$ws.Range("H65").Value = 4840.25
$ws.Range("I65").Value = 3227.5
$ws.Range("K65").Value = 16137.5
$ws.Range("M65").Value = -13017.5
$ws.Range("H68").Value = 118400
$ws.Range("J68").Value = 118400
$ws.Range("L68").Value = 118400
$ws.Range("N68").Value = -119898
$ws.Range("H69").Value = 29329.666
$ws.Range("I69").Value = 29329.666
$ws.Range("K69").Value = 29329.666
$ws.Range("M69").Value = -28580.666
$ws.Range("H71").Value = 118400
$ws.Range("J71").Value = 118400
$ws.Range("L71").Value = 355200
$ws.Range("N71").Value = -362688
$ws.Range("H72").Value = 29329.666
$ws.Range("I72").Value = 29329.666
$ws.Range("K72").Value = 87988.99800000001
$ws.Range("M72").Value = -84244.99800000001
$ws.Range("H74").Value = 96400
$ws.Range("J74").Value = 96400
$ws.Range("L74").Value = 96400
$ws.Range("N74").Value = -98148
$ws.Range("H77").Value = 96400
$ws.Range("J77").Value = 96400
$ws.Range("L77").Value = 289200
$ws.Range("N77").Value = -297936
$ws.Range("H96").Value = 31287
$ws.Range("J96").Value = 31287
$ws.Range("L96").Value = 31287
$ws.Range("N96").Value = -36779

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1779.4445
$ws.Range("I75").Value = 50
$ws.Range("J75").Value = 1995.625
$ws.Range("K75").Value = 150
$ws.Range("L75").Value = 5986.875
$ws.Range("M75").Value = 848
$ws.Range("N75").Value = -7982.875
$ws.Range("H78").Value = 1779.4445
$ws.Range("I78").Value = 50
$ws.Range("J78").Value = 1995.625
$ws.Range("K78").Value = 450
$ws.Range("L78").Value = 17960.625
$ws.Range("M78").Value = 4542
$ws.Range("N78").Value = -27944.625
$ws.Range("H80").Value = 1667
$ws.Range("J80").Value = 1999.5
$ws.Range("L80").Value = 5998.5
$ws.Range("N80").Value = -7870.5
$ws.Range("H83").Value = 1667
$ws.Range("J83").Value = 1999.5
$ws.Range("L83").Value = 17995.5
$ws.Range("N83").Value = -27355.5
$ws.Range("H86").Value = 286
$ws.Range("I86").Value = 292
$ws.Range("J86").Value = 280
$ws.Range("K86").Value = 876
$ws.Range("L86").Value = 840
$ws.Range("M86").Value = 310
$ws.Range("N86").Value = -3212
$ws.Range("H87").Value = 3004.25
$ws.Range("I87").Value = 3004.25
$ws.Range("K87").Value = 9012.75
$ws.Range("M87").Value = -7764.75
$ws.Range("H89").Value = 286
$ws.Range("I89").Value = 292
$ws.Range("J89").Value = 280
$ws.Range("K89").Value = 2628
$ws.Range("L89").Value = 2520
$ws.Range("M89").Value = 3300
$ws.Range("N89").Value = -14376
$ws.Range("H90").Value = 3004.25
$ws.Range("I90").Value = 3004.25
$ws.Range("K90").Value = 27038.25
$ws.Range("M90").Value = -20798.25
$ws.Range("H113").Value = 5292105
$ws.Range("J113").Value = 1325
$ws.Range("L113").Value = 3975
$ws.Range("N113").Value = -8315

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 48900
$ws.Range("J57").Value = 48900
$ws.Range("L57").Value = 48900
$ws.Range("N57").Value = -50540
$ws.Range("H110").Value = 40140.4
$ws.Range("J110").Value = 40140.4
$ws.Range("L110").Value = 40140.4
$ws.Range("N110").Value = -48320.4
$ws.Range("H126").Value = 3723.318
$ws.Range("J126").Value = 4027.5386
$ws.Range("L126").Value = 12082.6158
$ws.Range("N126").Value = -17022.6158
$ws.Range("H132").Value = 45870.73
$ws.Range("I132").Value = 7665.864
$ws.Range("K132").Value = 22997.592
$ws.Range("M132").Value = -20467.592
$ws.Range("H133").Value = 49000
$ws.Range("J133").Value = 49000
$ws.Range("L133").Value = 49000
$ws.Range("N133").Value = -59120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7642.3335
$ws.Range("I7").Value = 7514.3335
$ws.Range("K7").Value = 7514.3335
$ws.Range("M7").Value = -7402.3335
$ws.Range("H22").Value = 1171.4286
$ws.Range("J22").Value = 1080
$ws.Range("L22").Value = 1080
$ws.Range("N22").Value = -1670
$ws.Range("H27").Value = 1171.4286
$ws.Range("J27").Value = 1080
$ws.Range("L27").Value = 1080
$ws.Range("N27").Value = -1294
$ws.Range("H40").Value = 3086.611
$ws.Range("I40").Value = 2472.4375
$ws.Range("K40").Value = 2472.4375
$ws.Range("M40").Value = -2336.4375
$ws.Range("H126").Value = 7642.3335
$ws.Range("I126").Value = 7514.3335
$ws.Range("K126").Value = 22543.0005
$ws.Range("M126").Value = -20073.0005
$ws.Range("H132").Value = 7099.6924
$ws.Range("I132").Value = 5959.5186
$ws.Range("K132").Value = 17878.5558
$ws.Range("M132").Value = -15348.5558

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 958
$ws.Range("I100").Value = 933.3333
$ws.Range("J100").Value = 995
$ws.Range("K100").Value = 1866.6666
$ws.Range("L100").Value = 1990
$ws.Range("M100").Value = -1325.6666
$ws.Range("N100").Value = -3072
$ws.Range("I136").Value = 12363024
$ws.Range("J136").Value = 251893.12
$ws.Range("K136").Value = 37089072
$ws.Range("L136").Value = 755679.36
$ws.Range("M136").Value = -37086522
$ws.Range("N136").Value = -760779.36
